$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the study code typo (TLC01 -> TCL01) in the generated file name cells
$ws.Range("D2").Value = "TC01_ICDC_TCL01_Breed-GoldenRetriever_TSVData.xlsx"
$ws.Range("E2").Value = "TC01_ICDC_TCL01_Breed-GoldenRetriever_WebData.xlsx"

# Update the active selection to match the saved view state (D2 selected)
$ws.Range("D2").Select()
